$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24; existing rows 24-28 shift down to 25-29.
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with the new data record.
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 44855
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = 300000000
$ws.Cells.Item(24, 7).Value = "Espárragos"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 800
$ws.Cells.Item(24, 11).Value = 1000
$ws.Cells.Item(24, 12).Value = 1200
$ws.Cells.Item(24, 13).Value = 1100
$ws.Cells.Item(24, 14).Value = "$/kilo"
$ws.Cells.Item(24, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(24, 16).Value = 1100
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"
